$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 564.3333
$ws.Range("I28").Value = 568.9286
$ws.Range("K28").Value = 568.9286
$ws.Range("M28").Value = -83.92859999999996

$ws.Range("H48").Value = 5875
$ws.Range("J48").Value = 7666.6665
$ws.Range("L48").Value = 22999.9995
$ws.Range("N48").Value = -23583.9995

$ws.Range("H56").Value = 5875
$ws.Range("J56").Value = 7666.6665
$ws.Range("L56").Value = 22999.9995
$ws.Range("N56").Value = -24067.9995

$ws.Range("H103").Value = 2824.8
$ws.Range("I103").Value = 1840
$ws.Range("J103").Value = 3809.6
$ws.Range("K103").Value = 5520
$ws.Range("L103").Value = 11428.8
$ws.Range("M103").Value = -4934
$ws.Range("N103").Value = -12600.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1043.7778
$ws.Range("I122").Value = 1074.25
$ws.Range("K122").Value = 3222.75
$ws.Range("M122").Value = -772.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 238.25
$ws.Range("I12").Value = 317.33334
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 317.33334
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = -149.33334
$ws.Range("N12").Value = -337

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4899.8
$ws.Range("I7").Value = 5049.5
$ws.Range("J7").Value = 4800
$ws.Range("K7").Value = 5049.5
$ws.Range("L7").Value = 4800
$ws.Range("M7").Value = -4936.5
$ws.Range("N7").Value = -5026

$ws.Range("H132").Value = 5806.6
$ws.Range("I132").Value = 4505
$ws.Range("K132").Value = 13515
$ws.Range("M132").Value = -10985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 239.86667
$ws.Range("I2").Value = 57.166668
$ws.Range("J2").Value = 361.66666
$ws.Range("K2").Value = 57.166668
$ws.Range("L2").Value = 361.66666
$ws.Range("M2").Value = 55.833332
$ws.Range("N2").Value = -587.66666

$ws.Range("H3").Value = 5324.4
$ws.Range("J3").Value = 2872
$ws.Range("L3").Value = 2872
$ws.Range("N3").Value = -3104

$ws.Range("H4").Value = 637.8
$ws.Range("J4").Value = 496.66666
$ws.Range("L4").Value = 496.66666
$ws.Range("N4").Value = -720.66666

$ws.Range("H5").Value = 233
$ws.Range("I5").Value = 233
$ws.Range("K5").Value = 233
$ws.Range("M5").Value = -121

$ws.Range("H9").Value = 980
$ws.Range("J9").Value = 417.5
$ws.Range("L9").Value = 417.5
$ws.Range("N9").Value = -757.5

$ws.Range("H10").Value = 253400

$ws.Range("H11").Value = 850
$ws.Range("J11").Value = 950
$ws.Range("L11").Value = 950
$ws.Range("N11").Value = -1228

$ws.Range("H13").Value = 268.875
$ws.Range("I13").Value = 167.2
$ws.Range("J13").Value = 438.33334
$ws.Range("K13").Value = 167.2
$ws.Range("L13").Value = 438.33334
$ws.Range("M13").Value = -28.19999999999999
$ws.Range("N13").Value = -716.33334

$ws.Range("H14").Value = 165955.72
$ws.Range("J14").Value = 32158
$ws.Range("L14").Value = 32158
$ws.Range("N14").Value = -32494

$ws.Range("H15").Value = 25000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H17").Value = 556.7778
$ws.Range("I17").Value = 404.2
$ws.Range("J17").Value = 747.5
$ws.Range("K17").Value = 404.2
$ws.Range("L17").Value = 747.5
$ws.Range("M17").Value = -236.2
$ws.Range("N17").Value = -1083.5

$ws.Range("H70").Value = 300
$ws.Range("I70").Value = 300
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 300
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -30
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 300
$ws.Range("I73").Value = 300
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 300
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 636
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 7828
$ws.Range("I80").Value = 5949
$ws.Range("J80").Value = 10333.333
$ws.Range("K80").Value = 5949
$ws.Range("L80").Value = 10333.333
$ws.Range("M80").Value = -4951
$ws.Range("N80").Value = -12329.333

$ws.Range("H81").Value = 25000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H83").Value = 7828
$ws.Range("I83").Value = 5949
$ws.Range("J83").Value = 10333.333
$ws.Range("K83").Value = 29745
$ws.Range("L83").Value = 51666.665
$ws.Range("M83").Value = -24753
$ws.Range("N83").Value = -61650.665

$ws.Range("H84").Value = 25000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H132").Value = 1777.3846
$ws.Range("I132").Value = 1808.3636
$ws.Range("J132").Value = 1607
$ws.Range("K132").Value = 5425.0908
$ws.Range("L132").Value = 4821
$ws.Range("M132").Value = -2895.0908
$ws.Range("N132").Value = -9881

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4250
$ws.Range("I7").Value = 4250
$ws.Range("K7").Value = 4250
$ws.Range("M7").Value = -4138

$ws.Range("H18").Value = 250
$ws.Range("J18").Value = 250
$ws.Range("L18").Value = 250
$ws.Range("N18").Value = -594

$ws.Range("H46").Value = 3999.9333
$ws.Range("J46").Value = 4071.3572
$ws.Range("L46").Value = 4071.3572
$ws.Range("N46").Value = -4447.3572

$ws.Range("H55").Value = 556.4737
$ws.Range("I55").Value = 454
$ws.Range("J55").Value = 616.25
$ws.Range("K55").Value = 454
$ws.Range("L55").Value = 616.25
$ws.Range("M55").Value = -281
$ws.Range("N55").Value = -962.25

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1798

$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246

$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws.Range("H126").Value = 4250
$ws.Range("I126").Value = 4250
$ws.Range("K126").Value = 12750
$ws.Range("M126").Value = -10280

$ws.Range("H131").Value = 40000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 40000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 40000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -50080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2490
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws.Range("H132").Value = 1511.8334
$ws.Range("I132").Value = 1237.125
$ws.Range("K132").Value = 3711.375
$ws.Range("M132").Value = -1181.375

Write-Output "Edits applied successfully"